$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dataset (weekly Cilantro prices at Vega Monumental Concepcion)
# was refreshed: a new week of data was recorded, so the "Fecha" (date, column D)
# for every existing fortnight-pair of rows 154-187 shifts down by one pair of rows,
# carrying the "Origen" (column O) value along with it (only rows 160-163 actually
# change value, since elsewhere Origen stays "Region de Nuble" on both sides).
# The oldest pair (formerly rows 186-187) is preserved by appending it as new rows
# 188-189, and a brand-new date (44694) is recorded for the newest pair (154-155).

$ws.Cells.Item(154,4).Value = 44694
$ws.Cells.Item(155,4).Value = 44694
$ws.Cells.Item(156,4).Value = 44433
$ws.Cells.Item(157,4).Value = 44433
$ws.Cells.Item(158,4).Value = 44203
$ws.Cells.Item(159,4).Value = 44203
$ws.Cells.Item(160,4).Value = 44308
$ws.Cells.Item(160,15).Value = 'Región de Ñuble'
$ws.Cells.Item(161,4).Value = 44308
$ws.Cells.Item(161,15).Value = 'Región de Ñuble'
$ws.Cells.Item(162,4).Value = 44665
$ws.Cells.Item(162,15).Value = 'Provincia de Cautín'
$ws.Cells.Item(163,4).Value = 44665
$ws.Cells.Item(163,15).Value = 'Provincia de Cautín'
$ws.Cells.Item(164,4).Value = 44383
$ws.Cells.Item(165,4).Value = 44383
$ws.Cells.Item(166,4).Value = 44237
$ws.Cells.Item(167,4).Value = 44237
$ws.Cells.Item(168,4).Value = 44427
$ws.Cells.Item(169,4).Value = 44427
$ws.Cells.Item(170,4).Value = 44271
$ws.Cells.Item(171,4).Value = 44271
$ws.Cells.Item(172,4).Value = 44330
$ws.Cells.Item(173,4).Value = 44330
$ws.Cells.Item(174,4).Value = 44187
$ws.Cells.Item(175,4).Value = 44187
$ws.Cells.Item(176,4).Value = 44194
$ws.Cells.Item(177,4).Value = 44194
$ws.Cells.Item(178,4).Value = 44365
$ws.Cells.Item(179,4).Value = 44365
$ws.Cells.Item(180,4).Value = 44327
$ws.Cells.Item(181,4).Value = 44327
$ws.Cells.Item(182,4).Value = 44358
$ws.Cells.Item(183,4).Value = 44358
$ws.Cells.Item(184,4).Value = 44217
$ws.Cells.Item(185,4).Value = 44217
$ws.Cells.Item(186,4).Value = 44460
$ws.Cells.Item(187,4).Value = 44460

# --- Step 2: append the two rows that were pushed past the end of the range ---
# (identical to the old rows 186-187 content, which is what the data that
# used to be at rows 186-187 now looks like after the shift above).

# Row 188
$ws.Cells.Item(188,1).Value = 11
$ws.Cells.Item(188,2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(188,3).Value = 'Bíobío'
$ws.Cells.Item(188,4).Value = 44607
$ws.Cells.Item(188,5).Value = 8
$ws.Cells.Item(188,6).Value = 100112040
$ws.Cells.Item(188,7).Value = 'Cilantro'
$ws.Cells.Item(188,8).Value = 'Sin especificar'
$ws.Cells.Item(188,9).Value = 'Primera'
$ws.Cells.Item(188,10).Value = 200
$ws.Cells.Item(188,11).Value = 600
$ws.Cells.Item(188,12).Value = 700
$ws.Cells.Item(188,13).Value = 650
$ws.Cells.Item(188,14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(188,15).Value = 'Región de Ñuble'
$ws.Cells.Item(188,16).Value = 650
$ws.Cells.Item(188,17).Value = 1
$ws.Cells.Item(188,18).Value = 'Hortaliza'
$ws.Cells.Item(188,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 189
$ws.Cells.Item(189,1).Value = 11
$ws.Cells.Item(189,2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(189,3).Value = 'Bíobío'
$ws.Cells.Item(189,4).Value = 44607
$ws.Cells.Item(189,5).Value = 8
$ws.Cells.Item(189,6).Value = 100112040
$ws.Cells.Item(189,7).Value = 'Cilantro'
$ws.Cells.Item(189,8).Value = 'Sin especificar'
$ws.Cells.Item(189,9).Value = 'Segunda'
$ws.Cells.Item(189,10).Value = 100
$ws.Cells.Item(189,11).Value = 500
$ws.Cells.Item(189,12).Value = 500
$ws.Cells.Item(189,13).Value = 500
$ws.Cells.Item(189,14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(189,15).Value = 'Región de Ñuble'
$ws.Cells.Item(189,16).Value = 500
$ws.Cells.Item(189,17).Value = 1
$ws.Cells.Item(189,18).Value = 'Hortaliza'
$ws.Cells.Item(189,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

